# Capstone 1 report rev1
#
# Updates two pieces of body text in the deck:
#  - Slide 17 ("Modeling Results"): trim the parenthetical aside off the
#    Decision Tree mse bullet.
#  - Slide 4 ("Data wrangling"): extend the "Drop 60 lines" bullet with
#    more detail about which columns were affected.

$p = $ppt.ActivePresentation

# --- Slide 17: Decision Tree mse bullet -----------------------------------
$slide17 = $p.Slides.Item(17)
$body17 = $slide17.Shapes.Item(2).TextFrame.TextRange
$mseParagraph = $body17.Paragraphs(2, 1)
$mseRun = $mseParagraph.Runs(1, 1)
$mseRun.Text = "Decision Tree generated mse of 4.24 on held-out test data"

# --- Slide 4: Drop 60 lines bullet ----------------------------------------
$slide4 = $p.Slides.Item(4)
$body4 = $slide4.Shapes.Item(2).TextFrame.TextRange
$dropParagraph = $body4.Paragraphs(7, 1)
$dropRun = $dropParagraph.Runs(1, 1)
$dropRun.Text = "Drop 60 lines that have missing values across most of the skillset columns."
